$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.953.88'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '2.613.72'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.24'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.55'
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').Value = '2.613.09'
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('E10').Value = '  -3.67%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.19'
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.26'
$ws.Range('E14').Value = '  -3.09%  '
$ws.Range('D15').Value = '3.089.15'
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000180'
$ws.Range('E16').Value = '  -2.96%  '
$ws.Range('D17').Value = '66.819.81'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '2.602.32'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.67'
$ws.Range('E19').Value = '  -4.12%  '
$ws.Range('E20').Value = '  -4.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '355.19'
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.27'
$ws.Range('E22').Value = '  -3.29%  '
$ws.Range('E23').Value = '  -3.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.51'
$ws.Range('E24').Value = '  -4.69%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  -5.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '69.34'
$ws.Range('E27').Value = '  -2.93%  '
$ws.Range('D28').Value = '2.746.35'
$ws.Range('E28').Value = '  -1.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('E30').Value = '  -3.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '540.71'
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.17'
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('E33').Value = '  -4.18%  '
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('E35').Value = '  -1.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -5.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '158.40'
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.90'
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('E40').Value = '  -2.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.25'
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('E42').Value = '  -1.99%  '
$ws.Range('E43').Value = '  -3.80%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('E45').Value = '  -5.62%  '
$ws.Range('D46').Value = '0.0₆0294'
$ws.Range('E46').Value = '  -2.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.575'
$ws.Range('E47').Value = '  -3.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '150.58'
$ws.Range('E48').Value = '  -2.65%  '
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('E51').Value = '  -1.55%  '
